$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. Measures sheet: drop stray trailing formatted-only rows (53, 85, 98)
#    before inserting the new column so the used range collapses back to
#    A1:O52 (matches target dimension once the column insert makes it
#    A1:P52).
# -----------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("Measures")
$wsMeasures.Range("D53:E53").Clear()
$wsMeasures.Range("D85:E85").Clear()
$wsMeasures.Range("D98:E98").Clear()

# -----------------------------------------------------------------------
# 2. Insert a new column O ("item_num") ahead of the existing "comment"
#    column on every sheet that carries the Measures-style header row.
#    This shifts the old O ("comment") data to column P automatically.
# -----------------------------------------------------------------------
$sheetNames = @("Measures", "ID", "Dems", "Dates", "NewVars")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(15).Insert()
    $ws.Range("O1").Value = "item_num"
}

# -----------------------------------------------------------------------
# 3. Populate the new item_num column with 1 for every data row on the
#    Measures sheet (rows 2-52).
# -----------------------------------------------------------------------
$wsMeasures.Range("O2:O52").Value = 1

# -----------------------------------------------------------------------
# 4. Clean up the stray direct formatting that used to live in the J/M
#    columns on the Measures sheet.
#    - J2:J12 keep their numeric values but lose the applied number
#      format.
#    - J13:J18 / M19:M24 / M27:M28 were empty cells that existed only to
#      carry formatting - remove them outright.
# -----------------------------------------------------------------------
$wsMeasures.Range("J2:J12").ClearFormats()
$wsMeasures.Range("J13:J18").Clear()
$wsMeasures.Range("M19:M24").Clear()
$wsMeasures.Range("M27:M28").Clear()

# -----------------------------------------------------------------------
# 5. Defined names for the (hidden) filter databases need their ranges
#    widened by one column now that "comment" moved from O to P.
# -----------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Dems!_FilterDatabase") {
        $n.RefersTo = "=Dems!`$A`$1:`$P`$1"
    } elseif ($n.Name -eq "Measures!_FilterDatabase") {
        $n.RefersTo = "=Measures!`$A`$1:`$P`$52"
    }
}

# -----------------------------------------------------------------------
# 6. Sheet view / selection tidy-up to mirror the saved state in the
#    workbook: ID/Dems/Dates/NewVars all had their selection sitting on
#    the old "comment" column (I:I) which now lives at O:O, and the
#    active tab moves from NewVars to Measures, landing on O11:O52.
# -----------------------------------------------------------------------
$wb.Worksheets.Item("ID").Range("O1:O1048576").Select()
$wb.Worksheets.Item("Dems").Range("O1:O1048576").Select()
$wb.Worksheets.Item("Dates").Range("O1:O1048576").Select()
$wb.Worksheets.Item("NewVars").Range("O1:O1048576").Select()
$wsMeasures.Range("O11:O52").Select()
